$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 2, shifting existing data rows down one
$ws.Rows.Item(2).Insert()

# Pre-format the date-looking text columns as Text so Excel does not
# auto-convert them into date serial numbers
$ws.Range("A2").NumberFormat = "@"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("P2").NumberFormat = "@"

# Populate the new row 2 with the latest listing entry
$ws.Range("A2").Value = "2024-05-23"
$ws.Range("B2").Value = "노브랜드"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 168
$ws.Range("E2").Value = "삼성"
$ws.Range("F2").Value = 168
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 14000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2024-05-13"
$ws.Range("P2").Value = "2024-05-17"
$ws.Range("Q2").Value = 900000

# Reset the new row back to the plain/default style used by the rest of
# the data rows (Insert() otherwise inherits the bold/bordered header style)
$ws.Rows.Item(2).Style = "Normal"
